$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Relabel the "Enable" column header to "Status" (scheduleRun.enable row) ---
$ws.Range("C88").Value = "Status"

# --- 2. Insert two new rows before the old "scheduleRun.actions" row (old row 89) ---
#        for the new status toggle labels.
$ws.Rows("89:90").Insert()

$ws.Range("A89").Value = "scheduleRun.status.false"
$ws.Range("B89").Value = "Tắt"
$ws.Range("C89").Value = "Off"
$ws.Rows("89").RowHeight = 15

$ws.Range("A90").Value = "scheduleRun.status.true"
$ws.Range("B90").Value = "Bật"
$ws.Range("C90").Value = "On"
$ws.Rows("90").RowHeight = 15

# --- 3. Insert one new row after "scheduleRun.actions" (now row 91) for the new
#        confirmation-dialog message, before "button.close" (now row 92 -> 93).
$ws.Rows("92:92").Insert()

$ws.Range("A92").Value = "scheduleRun.message.modifyStatus"
$ws.Range("B92").Value = "Bạn có thật sự muốn thay đổi trạng thái?"
$ws.Range("C92").Value = "Do you want to change status?"
$ws.Rows("92").RowHeight = 15

# --- 4. The row-insert above propagated the old, now-stale "A89" cell format
#        (a leftover one-off font variant) onto A91/A92. Re-normalize those two
#        cells back to the plain style used by the rest of the key column so
#        that stray format doesn't linger on the new layout.
$ws.Range("A10").Copy()
$ws.Range("A91").PasteSpecial(-4122)
$ws.Range("A92").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. Move the selection to reflect the author's last-active cell ---
$ws.Range("C92").Select() | Out-Null
